# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Adds a new worker (JEFFERSON JOSE CONEO RODRIGUEZ, doc 1050960720) with
# period 2509, gives the existing worker KEITNER MARTINEZ BARRIOS (doc
# 1047475681) an additional overdue period 2509, and refreshes the
# document's summary totals (Valor Mora, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for two new data rows right after the current last data
#     row (21 is the last "regular" styled row, 22 is the closing row with
#     the heavier bottom border). Insert two rows before the closing row and
#     restyle them like a normal data row (copy format+values from row 21),
#     so the closing-row look stays at the bottom of the table.
$ws.Rows("22:23").Insert()
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))
$ws.Range("B21:J21").Copy($ws.Range("B23:J23"))

# Row 22 keeps being KEITNER MARTINEZ BARRIOS, period 2508 (same as before,
# just no longer the last row of the table)
$ws.Range("E22").Value = "2508"

# Row 23: brand-new worker, period 2509
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1050960720"
$ws.Range("D23").Value = "JEFFERSON JOSE CONEO RODRIGUEZ"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

# Row 24 (old closing row, pushed down automatically by the insert above):
# same worker as row 22 (KEITNER MARTINEZ BARRIOS) but for the new period
# 2509 - this becomes the new closing/last row of the table.
$ws.Range("E24").Value = "2509"

# --- Column D ("Nombre Trabajador") needs to fit the new, longer name.
$ws.Columns("D:D").AutoFit()

# --- Refresh the summary block above the table.
$ws.Range("E11").Value = 478241   # VALOR MORA total
$ws.Range("C13").Value = 5        # Cant. Trabajadores
$ws.Range("F13").Value = 8        # Cant. Periodos
